$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new test data: "AssertValidLogin" assertion for the "LoginWithValidData" row
$ws.Range("D7").Value = "AssertValidLogin"

# Apply a thin box border around every cell of the data table (header + 6 data rows)
$tableRange = $ws.Range("A1:D7")
$tableRange.Borders.LineStyle = 1
$tableRange.Borders.Weight = 2

# Clear the fill of the two blank rows below the table (previously greenish fill on col A)
$blankRange = $ws.Range("A8:D9")
$blankRange.Interior.Pattern = 0

# Move the active selection to C11
$ws.Range("C11").Select()
